$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-05-10 Friday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-05-11 Saturday", 2) | Out-Null

# Update each answer cell in the table by its (row, column) position so that
# cells whose old/new text overlap with other cells values do not collide.
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "31÷6=5, 1"

$cell = $t.Cell(1, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "83÷4=20, 3"

$cell = $t.Cell(1, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "65÷3=21, 2"

$cell = $t.Cell(1, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "72÷5=14, 2"

$cell = $t.Cell(1, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "13÷8=1, 5"

$cell = $t.Cell(5, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "61÷7=8, 5"

$cell = $t.Cell(5, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "79÷4=19, 3"

$cell = $t.Cell(5, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "83÷5=16, 3"

$cell = $t.Cell(5, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "46÷7=6, 4"

$cell = $t.Cell(5, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "28÷5=5, 3"

$cell = $t.Cell(9, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "66÷4=16, 2"

$cell = $t.Cell(9, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "92÷3=30, 2"

$cell = $t.Cell(9, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "58÷8=7, 2"

$cell = $t.Cell(9, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "73÷4=18, 1"

$cell = $t.Cell(9, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "83÷9=9, 2"

$cell = $t.Cell(13, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "31÷3=10, 1"

$cell = $t.Cell(13, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "44÷7=6, 2"

$cell = $t.Cell(13, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "64÷5=12, 4"

$cell = $t.Cell(13, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "57÷6=9, 3"

$cell = $t.Cell(13, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "15÷2=7, 1"

$cell = $t.Cell(17, 1)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "66÷3=22, 0"

$cell = $t.Cell(17, 2)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "16÷9=1, 7"

$cell = $t.Cell(17, 3)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "87÷3=29, 0"

$cell = $t.Cell(17, 4)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "21÷8=2, 5"

$cell = $t.Cell(17, 5)
$rng = $cell.Range
$rng.MoveEnd(1, -1) | Out-Null
$rng.Text = "40÷3=13, 1"
